# ---------------------------------------------------------------------------
# "updated part 2 results" - rewrite the status-code breakdown table so it
# reflects the new counts (rows 2-26, columns A:C). Row 1 (header) is
# untouched. New rows are appended below the previous last row (row 7) and
# the sheet dimension grows from A1:C7 to A1:C26 automatically as values are
# written.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "status code" labels and is formatted (bold, bordered,
# centered) the same way throughout the table - that formatting already sits
# on A2 in the original sheet, so we clone it for every row we touch/add.
# A second helper clone is kept in Text format, because most status codes
# ("200", "404", ...) look like numbers and must be forced to stay text -
# exactly like the source file already stores them.
$plainTemplate = $ws.Range("ZZ1")
$textTemplate = $ws.Range("ZZ2")
$ws.Range("A2").Copy($plainTemplate)
$ws.Range("A2").Copy($textTemplate)
$textTemplate.NumberFormat = "@"
$textTemplate.Value = "x"

# Row 2: status code '200'
$textTemplate.Copy($ws.Range("A2"))
$ws.Range("A2").Value = "200"
$ws.Range("B2").Value = 3198
$ws.Range("C2").Value = 6232

# Row 3: status code '202'
$textTemplate.Copy($ws.Range("A3"))
$ws.Range("A3").Value = "202"
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = $null

# Row 4: status code '400'
$textTemplate.Copy($ws.Range("A4"))
$ws.Range("A4").Value = "400"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4

# Row 5: status code '401'
$textTemplate.Copy($ws.Range("A5"))
$ws.Range("A5").Value = "401"
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = 1

# Row 6: status code '403'
$textTemplate.Copy($ws.Range("A6"))
$ws.Range("A6").Value = "403"
$ws.Range("B6").Value = 1478
$ws.Range("C6").Value = 2879

# Row 7: status code '404'
$textTemplate.Copy($ws.Range("A7"))
$ws.Range("A7").Value = "404"
$ws.Range("B7").Value = 104
$ws.Range("C7").Value = 74

# Row 8: status code '405'
$textTemplate.Copy($ws.Range("A8"))
$ws.Range("A8").Value = "405"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 5

# Row 9: status code '406'
$textTemplate.Copy($ws.Range("A9"))
$ws.Range("A9").Value = "406"
$ws.Range("B9").Value = $null
$ws.Range("C9").Value = 82

# Row 10: status code '409'
$textTemplate.Copy($ws.Range("A10"))
$ws.Range("A10").Value = "409"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = $null

# Row 11: status code '410'
$textTemplate.Copy($ws.Range("A11"))
$ws.Range("A11").Value = "410"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 5

# Row 12: status code '415'
$textTemplate.Copy($ws.Range("A12"))
$ws.Range("A12").Value = "415"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = $null

# Row 13: status code '429'
$textTemplate.Copy($ws.Range("A13"))
$ws.Range("A13").Value = "429"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = $null

# Row 14: status code '447'
$textTemplate.Copy($ws.Range("A14"))
$ws.Range("A14").Value = "447"
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = 6

# Row 15: status code '468'
$textTemplate.Copy($ws.Range("A15"))
$ws.Range("A15").Value = "468"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = $null

# Row 16: status code '500'
$textTemplate.Copy($ws.Range("A16"))
$ws.Range("A16").Value = "500"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = 17

# Row 17: status code '502'
$textTemplate.Copy($ws.Range("A17"))
$ws.Range("A17").Value = "502"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 4

# Row 18: status code '503'
$textTemplate.Copy($ws.Range("A18"))
$ws.Range("A18").Value = "503"
$ws.Range("B18").Value = 8
$ws.Range("C18").Value = $null

# Row 19: status code '504'
$textTemplate.Copy($ws.Range("A19"))
$ws.Range("A19").Value = "504"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = $null

# Row 20: status code '520'
$textTemplate.Copy($ws.Range("A20"))
$ws.Range("A20").Value = "520"
$ws.Range("B20").Value = $null
$ws.Range("C20").Value = 4

# Row 21: status code '523'
$textTemplate.Copy($ws.Range("A21"))
$ws.Range("A21").Value = "523"
$ws.Range("B21").Value = $null
$ws.Range("C21").Value = 2

# Row 22: status code '525'
$textTemplate.Copy($ws.Range("A22"))
$ws.Range("A22").Value = "525"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = $null

# Row 23: status code '526'
$textTemplate.Copy($ws.Range("A23"))
$ws.Range("A23").Value = "526"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = $null

# Row 24: status code '530'
$textTemplate.Copy($ws.Range("A24"))
$ws.Range("A24").Value = "530"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = $null

# Row 25: status code '555'
$textTemplate.Copy($ws.Range("A25"))
$ws.Range("A25").Value = "555"
$ws.Range("B25").Value = $null
$ws.Range("C25").Value = 1

# Row 26: status code 'Attempt failed'
$plainTemplate.Copy($ws.Range("A26"))
$ws.Range("A26").Value = "Attempt failed"
$ws.Range("B26").Value = 5159
$ws.Range("C26").Value = 678

# Drop the scratch column used only to carry the text-number-format style -
# keeps the sheet dimension at A1:C26 with no stray data outside of it.
$ws.Range("ZZ1").EntireColumn.Delete()
